# Apply cryptos list update (prices / 1h volume % changes, and a few
# coin re-rankings in rows 46-51) per the authoritative diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '68.371.82'
$ws.Range("E2").Value = '  -1.55%  '

# Row 3
$ws.Range("D3").Value = '3.842.63'
$ws.Range("E3").Value = '  -1.39%  '

# Row 4
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.01%  '

# Row 5
$ws.Range("D5").Value = '''601.79'
$ws.Range("E5").Value = '  -0.40%  '

# Row 6
$ws.Range("D6").Value = '''168.83'
$ws.Range("E6").Value = '  -0.63%  '

# Row 7
$ws.Range("D7").Value = '3.840.62'
$ws.Range("E7").Value = '  -1.40%  '

# Row 8
$ws.Range("E8").Value = '  +0.06%  '

# Row 9
$ws.Range("E9").Value = '  -1.03%  '

# Row 10
$ws.Range("E10").Value = '  -2.41%  '

# Row 11
$ws.Range("D11").Value = '''6.47'
$ws.Range("E11").Value = '  +1.06%  '

# Row 12
$ws.Range("E12").Value = '  -2.35%  '

# Row 13
$ws.Range("E13").Value = '  +4.49%  '

# Row 14
$ws.Range("D14").Value = '''37.03'
$ws.Range("E14").Value = '  -3.19%  '

# Row 15
$ws.Range("D15").Value = '4.492.72'
$ws.Range("E15").Value = '  -0.33%  '

# Row 16
$ws.Range("D16").Value = '3.840.05'
$ws.Range("E16").Value = '  -1.42%  '

# Row 17
$ws.Range("D17").Value = '68.460.84'
$ws.Range("E17").Value = '  -1.45%  '

# Row 18
$ws.Range("D18").Value = '''18.54'
$ws.Range("E18").Value = '  -1.29%  '

# Row 19
$ws.Range("D19").Value = '''7.37'
$ws.Range("E19").Value = '  -3.40%  '

# Row 20
$ws.Range("E20").Value = '  -1.17%  '

# Row 21
$ws.Range("D21").Value = '''11.13'
$ws.Range("E21").Value = '  +0.56%  '

# Row 22
$ws.Range("D22").Value = '''470.46'
$ws.Range("E22").Value = '  -3.93%  '

# Row 23
$ws.Range("D23").Value = '''0.732'
$ws.Range("E23").Value = '  -1.56%  '

# Row 24
$ws.Range("E24").Value = '  -3.93%  '

# Row 25
$ws.Range("D25").Value = '''83.46'
$ws.Range("E25").Value = '  -2.16%  '

# Row 26
$ws.Range("E26").Value = '  -2.85%  '

# Row 27
$ws.Range("D27").Value = '''12.10'
$ws.Range("E27").Value = '  -2.41%  '

# Row 28
$ws.Range("D28").Value = '''10.24'
$ws.Range("E28").Value = '  +0.77%  '

# Row 29
$ws.Range("E29").Value = '  -0.02%  '

# Row 30
$ws.Range("E30").Value = '  -0.95%  '

# Row 31
$ws.Range("D31").Value = '3.993.41'
$ws.Range("E31").Value = '  -1.30%  '

# Row 32
$ws.Range("E32").Value = '  -1.40%  '

# Row 33
$ws.Range("D33").Value = '''31.42'
$ws.Range("E33").Value = '  -1.72%  '

# Row 34
$ws.Range("E34").Value = '  -3.65%  '

# Row 35
$ws.Range("E35").Value = '  -2.98%  '

# Row 36
$ws.Range("D36").Value = '3.809.89'
$ws.Range("E36").Value = '  -1.50%  '

# Row 37
$ws.Range("D37").Value = '''3.79'
$ws.Range("E37").Value = '  +11.48%  '

# Row 38
$ws.Range("E38").Value = '  -2.25%  '

# Row 39
$ws.Range("E39").Value = '  -1.52%  '

# Row 40
$ws.Range("D40").Value = '''0.139'
$ws.Range("E40").Value = '  -2.64%  '

# Row 41
$ws.Range("E41").Value = '  -2.87%  '

# Row 42
$ws.Range("E42").Value = '  +0.10%  '

# Row 43
$ws.Range("E43").Value = '  -3.33%  '

# Row 44
$ws.Range("E44").Value = '  -4.37%  '

# Row 45
$ws.Range("D45").Value = '''8.73'
$ws.Range("E45").Value = '  +0.37%  '

# Row 46
$ws.Range("B46").Value = 'FLOKI'
$ws.Range("C46").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D46").Value = '''0.000295'
$ws.Range("E46").Value = '  +6.36%  '

# Row 47
$ws.Range("B47").Value = 'Bittensor'
$ws.Range("C47").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D47").Value = '''419.18'
$ws.Range("E47").Value = '  -4.22%  '

# Row 48
$ws.Range("B48").Value = 'USDe'
$ws.Range("C48").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D48").Value = '''1.00'
$ws.Range("E48").Value = '  -0.01%  '

# Row 49
$ws.Range("D49").Value = '''46.98'
$ws.Range("E49").Value = '  -2.22%  '

# Row 50
$ws.Range("D50").Value = '''141.76'
$ws.Range("E50").Value = '  +0.23%  '

# Row 51
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '''26.05'
$ws.Range("E51").Value = '  +3.61%  '
